# Auto-generated edit script applying cell-level changes from the commit diff.
# Numeric-looking values in columns D and E are entered with a leading single
# quote (text-prefix) so Excel stores them as literal text, exactly like the
# original workbook (which stored all values as inline/shared strings), instead
# of silently converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.28"
$ws.Range("E2").Value = "'-0.06%"
$ws.Range("D3").Value = "'41.05"
$ws.Range("E3").Value = "'0.54%"
$ws.Range("D4").Value = "'5.198"
$ws.Range("E4").Value = "'1.65%"
$ws.Range("D5").Value = "'0.07676"
$ws.Range("E5").Value = "'0.67%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.630"
$ws.Range("E6").Value = "'1.64%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9142"
$ws.Range("E7").Value = "'1.25%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.430"
$ws.Range("E8").Value = "'0.16%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1226"
$ws.Range("E9").Value = "'10.32%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1821"
$ws.Range("E10").Value = "'2.00%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09111"
$ws.Range("E11").Value = "'-0.28%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04252"
$ws.Range("E12").Value = "'1.19%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1052"
$ws.Range("E13").Value = "'-0.04%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001257"
$ws.Range("E14").Value = "'0.63%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005751"
$ws.Range("E15").Value = "'1.43%"
$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").Value = "'0.007509"
$ws.Range("E16").Value = "'1,904.49%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.342"
$ws.Range("E17").Value = "'-0.17%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.298"
$ws.Range("E18").Value = "'1.30%"
$ws.Range("D20").Value = "'7.393"
$ws.Range("E20").Value = "'12.32%"
$ws.Range("D21").Value = "'0.1382"
$ws.Range("E21").Value = "'1.29%"
$ws.Range("D22").Value = "'0.2713"
$ws.Range("E22").Value = "'-4.13%"
$ws.Range("D23").Value = "'0.04014"
$ws.Range("E23").Value = "'-1.58%"
$ws.Range("D24").Value = "'0.001263"
$ws.Range("E24").Value = "'2.67%"
$ws.Range("D25").Value = "'0.004377"
$ws.Range("E25").Value = "'6.15%"
$ws.Range("E26").Value = "'-0.01%"
$ws.Range("D38").Value = "'0.02501"
$ws.Range("E38").Value = "'3.61%"
$ws.Range("D39").Value = "'0.05303"
$ws.Range("E39").Value = "'2.22%"
$ws.Range("D40").Value = "'0.007836"
$ws.Range("E40").Value = "'0.87%"
$ws.Range("D41").Value = "'0.1314"
$ws.Range("E41").Value = "'0.88%"
$ws.Range("D42").Value = "'0.006591"
$ws.Range("E42").Value = "'-6.50%"
$ws.Range("E43").Value = "'-4.63%"
$ws.Range("D44").Value = "'0.008026"
$ws.Range("E44").Value = "'-8.61%"
$ws.Range("D45").Value = "'0.3038"
$ws.Range("E45").Value = "'-8.90%"
$ws.Range("D46").Value = "'0.00006714"
$ws.Range("E46").Value = "'-3.20%"
$ws.Range("E47").Value = "'0.04%"
$ws.Range("D48").Value = "'0.3501"
$ws.Range("E48").Value = "'1,008.11%"
$ws.Range("D49").Value = "'0.003103"
$ws.Range("E49").Value = "'-26.13%"
$ws.Range("E50").Value = "'0.04%"
$ws.Range("E51").Value = "'0.04%"
